
# myRIO Demo Board BOM -- Rev 1.1 update
#
# - Silkscreen/ref-des renames:
#     C3 (myRIO MXP connector) : "myRIO MXP" -> "MXP1"
#     C7 (7-seg display)       : "7SEG1"     -> "SEV_SEG1"
# - Package cleanup for the 7-seg display:
#     H7 : "10-DIP (0.600"", 15.24mm)" -> "10-DIP"
# - Type column cleanup for the myRIO MXP connector (now matches the
#   plain "Through Hole" used by the rest of the BOM rows):
#     I3 : "Through Hole, Right Angle" -> "Through Hole"
# - Selection moved to D10.
#
# Note: writes are ordered so that newly-introduced shared strings land in
# the workbook in the same order the reference edit produced them
# (10-DIP, SEV_SEG1, MXP1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H7").Value = "10-DIP"
$ws.Range("C7").Value = "SEV_SEG1"
$ws.Range("C3").Value = "MXP1"
$ws.Range("I3").Value = "Through Hole"

$ws.Range("D10").Select()
